$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.303.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.327.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.59%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.20%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  +2.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.355.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("E10").Value = "  +7.45%  "

$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.28%  "

$ws.Range("E13").Value = "  +1.84%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.84%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.777.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.166.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.358.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("E20").Value = "  +3.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.82%  "

$ws.Range("E22").Value = "  +5.79%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("E25").Value = "  +7.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("E27").Value = "  +5.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("E29").Value = "  +9.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0740"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.19%  "

$ws.Range("E31").Value = "  +4.15%  "

$ws.Range("E32").Value = "  +3.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.959"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.38%  "

$ws.Range("E38").Value = "  +7.14%  "

$ws.Range("E39").Value = "  +4.02%  "

$ws.Range("E40").Value = "  +7.05%  "

$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("E46").Value = "  +3.29%  "

$ws.Range("E47").Value = "  +3.55%  "

$ws.Range("E48").Value = "  +2.34%  "

$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("E50").Value = "  +4.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.43%  "
